# Update 13C-MFA files (run and result) for SC and IO under WT-batch and chemostats
# - FluxData: insert a new flux row (EX_glc__D_e.f) after BIOMASS.f, refresh
#   values for the rows whose data changed, and append a duplicated
#   DIL_ade_d1.f row at the bottom.
# - Refresh view state (zoom, selection, active sheet) on all three sheets;
#   FluxData becomes the active / selected tab instead of MSData.

$wb = $excel.ActiveWorkbook

$msData = $wb.Worksheets.Item("MSData")
$fluxData = $wb.Worksheets.Item("FluxData")
$tracerData = $wb.Worksheets.Item("TracerData")

# --- FluxData: insert the new EX_glc__D_e.f flux row at row 3 ---------------
$fluxData.Rows.Item(3).Insert()

$fluxData.Range("A3").Value = "EX_glc__D_e.f"
$fluxData.Range("B3").Value = 1.6916514664188
$fluxData.Range("C3").Value = 0.354858945709085
$fluxData.Rows.Item(3).RowHeight = 13.8

# BIOMASS.f row: only the basis (C2) value was refreshed
$fluxData.Range("C2").Value = 0.0001

# EX_c5sugal_e.f row (pushed down to row 4) got refreshed flux values
$fluxData.Range("B4").Value = 0.115740740740741
$fluxData.Range("C4").Value = 0.039890770682841

# Last DIL_* row (DIL_val__L_d1.f, now row 36) keeps ht=15
$fluxData.Rows.Item(36).RowHeight = 15

# New bottom row 37: duplicate of the former last row (DIL_ade_d1.f)
$fluxData.Range("A37").Value = "DIL_ade_d1.f"
$fluxData.Range("B37").Value = 100
$fluxData.Range("C37").Value = 0.0001
$fluxData.Rows.Item(37).RowHeight = 13.8

# --- View state: zoom 55 -> 95 on every sheet --------------------------------
$msData.Activate()
$msData.Range("A2").Select() | Out-Null
$excel.ActiveWindow.Zoom = 95

$fluxData.Activate()
$fluxData.Range("A1").Select() | Out-Null
$excel.ActiveWindow.Zoom = 95

$tracerData.Activate()
$tracerData.Range("A1").Select() | Out-Null
$excel.ActiveWindow.Zoom = 95

# FluxData becomes the active / selected tab (was MSData)
$fluxData.Activate()
